$data = @(
  @("FAPs","Wnt5a","Mcam","ECs",3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,40.23443200000001,120.703296,0.5194057602668869,0.5194057602668869,419.8158538664534,3778.34268479808,0.5043997112265909,0.504399711226591),
  @("FAPs","Wnt5a","Mcam","FAPs",3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,2.585148666666667,7.755446,0.03337293561427507,0.03337293561427508,26.97407024084222,242.76663216758,0.03240876473525146,0.03240876473525147),
  @("FAPs","Wnt5a","Mcam","M2",3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,0.525608,1.576824,0.006785328120013172,0.006785328120013173,5.484321769946666,49.35889592952,0.006589294547972889,0.00658929454797289),
  @("FAPs","Wnt5a","Mcam","sCs",3,1,10.43424333333333,31.30273,0.9711091978791583,0.9711091978791584,3,1,34.11724066666667,102.351722,0.4404359759988248,0.4404359759988249,355.9875909778955,3203.88831880106,0.4277114273693429,0.4277114273693431),
  @("sCs","Wnt5a","Mcam","ECs",2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,40.23443200000001,120.703296,0.5194057602668869,0.5194057602668869,12.489652850304,112.406875652736,0.01500604904029592,0.01500604904029593),
  @("sCs","Wnt5a","Mcam","FAPs",2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,2.585148666666667,7.755446,0.03337293561427507,0.03337293561427508,0.8024870194039999,7.222383174636,0.0009641708790236088,0.0009641708790236091),
  @("sCs","Wnt5a","Mcam","M2",2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,0.525608,1.576824,0.006785328120013172,0.006785328120013173,0.163160286576,1.468442579184,0.0001960335720402828,0.0001960335720402828),
  @("sCs","Wnt5a","Mcam","sCs",2,0.6666666666666666,0.310422,0.9312659999999999,0.02889080212084161,0.02889080212084161,3,1,34.11724066666667,102.351722,0.4404359759988248,0.4404359759988249,10.590742082228,95.316678740052,0.01272454862948179,0.0127245486294818)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
